$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 32: fill in the previously-missing "tested" value ---
$ws.Range("B32").Value = 827

# --- Row 33: new data for Apr 6, 2020 (date serial 43927) ---
$ws.Range("A33").Value = 43927
$ws.Range("C33").Value = 1404
$ws.Range("D33").Value = 10
$ws.Range("F33").Value = 15
$ws.Range("H33").Value = 235
$ws.Range("J33").Value = 300
$ws.Range("L33").Value = 254
$ws.Range("N33").Value = 235
$ws.Range("P33").Value = 172
$ws.Range("R33").Value = 109
$ws.Range("T33").Value = 71
$ws.Range("V33").Value = 3
$ws.Range("X33").Value = 657
$ws.Range("Y33").Value = 739
$ws.Range("Z33").Value = 8
$ws.Range("AA33").Value = 269
$ws.Range("AB33").Value = 102
$ws.Range("AC33").Value = 19
$ws.Range("AD33").Value = 42
$ws.Range("AE33").Value = 104
$ws.Range("AF33").Value = 4
$ws.Range("AG33").Value = 8
$ws.Range("AH33").Value = 74
$ws.Range("AI33").Value = 30
$ws.Range("AJ33").Value = 28
$ws.Range("AK33").Value = 7
$ws.Range("AL33").Value = 23
$ws.Range("AM33").Value = 13
$ws.Range("AN33").Value = 25
$ws.Range("AO33").Value = 31
$ws.Range("AP33").Value = 14
$ws.Range("AQ33").Value = 713
$ws.Range("AR33").Value = 19
$ws.Range("AS33").Value = 16
$ws.Range("AT33").Value = 5
$ws.Range("AU33").Value = 21
$ws.Range("AV33").Value = 1
$ws.Range("AW33").Value = 8
$ws.Range("AX33").Value = 1
$ws.Range("AY33").Value = 1
$ws.Range("AZ33").Value = 6
$ws.Range("BA33").Value = 2
$ws.Range("BB33").Value = 12
$ws.Range("BC33").Value = 2
$ws.Range("BD33").Value = 8
$ws.Range("BE33").Value = 14
$ws.Range("BG33").Value = 34
$ws.Range("BH33").Value = 3
$ws.Range("BI33").Value = 63

# --- View state: selection + scroll position ---
$ws.Range("BK34").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 51
$excel.ActiveWindow.ScrollRow = 1
